$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.183.92"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "2.445.54"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.21"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.18"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "2.440.65"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  +2.77%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.45"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "2.868.33"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "62.083.35"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "2.434.69"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.75"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.49"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  -5.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.68"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.16"
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "604.44"
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("D28").Value = "0.0₃0967"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "2.568.59"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.41"
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.89"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.135"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.89"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.49"
$ws.Range("E39").Value = "  +5.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.41"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.28"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.16"
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.02"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.63"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "0.0₆0265"
$ws.Range("E48").Value = "  +19.07%  "
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.84"
$ws.Range("E51").Value = "  +1.08%  "
